# Add a "p-value" column (column C) to worksheets tbl5-tbl8, with the
# header styled like the existing headers (bold + centered) and numeric
# p-values for each data row.

$wb = $excel.ActiveWorkbook

# sheet name -> ordered list of p-values for rows 2..N
$pvalues = @{
    "tbl5" = @(0.79, 0.5600000000000001, 0.07000000000000001, 0.18, 0.08, 0.41, 0.14, 0.67, 0.95, 0.06)
    "tbl6" = @(0.06, 0.34, 0.6, 0.32, 0.32, 0.54, 0.64, 0.86, 0.24, 0.62, 0.19, 0.03, 0.9399999999999999, 0.42, 0.13, 0.76, 0.59, 0.21, 0.25, 0.62)
    "tbl7" = @(0.05, 0.44, 0.73, 0.43, 0.03, 0.95, 0.31, 0.08, 0.79, 0.5600000000000001, 0.22, 0.25, 0.61)
    "tbl8" = @(0.05, 0.6, 0.39, 0.04, 0.03, 0.9399999999999999, 0.29, 0.08, 0.82, 0.5600000000000001, 0.24, 0.31, 0.61)
}

foreach ($sheetName in $pvalues.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Header cell: same look as the existing header row (bold + centered),
    # matching the style already used by A1/B1.
    $ws.Range("C1").Value = "p-value"
    $ws.Range("C1").Font.Bold = $true
    $ws.Range("C1").HorizontalAlignment = -4108

    $rows = $pvalues[$sheetName]
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $rowNum = $i + 2
        $ws.Cells.Item($rowNum, 3).Value = $rows[$i]
    }
}
